$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A8").Value = "rStorageNVRAM"
$ws.Range("B8").Value = "ctor_dtor_leak"
$ws.Range("C8").Value = "/home/rdkv-core/cov/cov-analysis-linux64-2023.6.0/bin/device/StorageManager/refactored/mgr/rdkStorageNVRAM.cpp"
$ws.Range("D8").Value = 36
$ws.Range("C14").Select()
